# Split the sentence ". Utiliza el componente en el archivo de entrada de
# ViteJs (main.tsx)" into several runs and append a new sentence about
# rendering the UserProfile component, per the "feat: add react-query
# example" commit.
#
# Target run structure (all runs share <w:lang w:val="es-EC"/>):
#   1) ". Utiliza el componente en el archivo de entrada de ViteJs ("
#   2) "App"
#   3) ".tsx)"
#   4) " para renderizar el nuevo componente UserProfile dentro de la
#        carpeta \u201ccomponents\u201d."   (xml:space="preserve")

$d = $word.ActiveDocument

# --- Step 1: narrow in on "main.tsx)" (unique in the document) and turn it
# into "App.tsx)". Starting the replace mid-run (not at the run boundary
# right after the preceding hyperlink) keeps the original run's own
# formatting (es-EC) instead of inheriting the hyperlink's rStyle.
$r = $d.Content
$r.Find.Execute("main.tsx)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "App.tsx)"

# --- Step 2: split "App" into its own run by toggling a direct-formatting
# property on and back off; Word (and this host) splits runs at the
# boundaries of the formatted sub-range while preserving the paragraph
# language formatting on every resulting run.
$appRun = $d.Range($r.Start, $r.Start + 3)
$appRun.Font.Bold = 1
$appRun.Font.Bold = 0

# --- Step 3: append the new trailing sentence after ".tsx)".
$newSentence = " para renderizar el nuevo componente UserProfile dentro de la carpeta " + [char]0x201C + "components" + [char]0x201D + "."
$r.InsertAfter($newSentence)

# --- Step 4: the run just inserted via InsertAfter has no rPr yet; locate
# it again with Find (re-resolving through Find makes the LanguageID
# assignment stick) and stamp it with the same es-EC language as its
# neighbours so it serializes with <w:rPr><w:lang w:val="es-EC"/></w:rPr>.
$newRun = $d.Content
$newRun.Find.Execute("para renderizar el nuevo componente UserProfile dentro de la carpeta", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newRun.LanguageID = "es-EC"
